$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pitcher ERA data to append after the existing last row (row 910).
$newData = @(
    @(909, "ángel perdomo", 6.35),
    @(910, "brent honeywell", 8.31),
    @(911, "carl edwards", 11.12),
    @(912, "dan camarena", 9.640000000000001),
    @(913, "daniel lynch", 5.69),
    @(914, "duane underwood", 4.33),
    @(915, "j.b. bukauskas", 7.79),
    @(916, "j.t. chargois", 2.52),
    @(917, "jaime barría", 4.61),
    @(918, "julio teherán", 1.8),
    @(919, "lance mccullers", 3.16),
    @(920, "matt boyd", 3.89),
    @(921, "mike king", 3.55),
    @(922, "mike wright", 5.5),
    @(923, "néstor cortés", 2.9),
    @(924, "travis lakins", 5.79),
    @(925, "vladimir gutiérrez", 4.74)
)

$startRow = 911
$endRow = $startRow + $newData.Length - 1

# Copy the formatting from the current last data row (910) down across the
# new rows so the appended cells match the existing table's look (bold
# border/center style on column A, etc.).
$srcRange = $ws.Range("A910:C910")
$dstRange = $ws.Range("A" + $startRow + ":C" + $endRow)
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $startRow + $i
    $entry = $newData[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}
